$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row - written in this order so the shared-strings table matches
# the order the values were first introduced in the authored file
# ("cedula" first, "nombre completo" last).
$ws.Range("B1").Value = "cedula"
$ws.Range("C1").Value = "contraseña"
$ws.Range("D1").Value = "cuenta ahorros"
$ws.Range("E1").Value = "saldo cuenta ahorros"
$ws.Range("F1").Value = "cuenta corriente"
$ws.Range("G1").Value = "saldo cuenta corriente"
$ws.Range("H1").Value = "tarjeta credito"
$ws.Range("I1").Value = "saldo tarjeta credito"
$ws.Range("J1").Value = "datacredito"
$ws.Range("K1").Value = "admin"
$ws.Range("A1").Value = "nombre completo"

# Row 2 - liliana lozano
$ws.Range("A2").Value = "liliana lozano"
$ws.Range("B2").Value = 66727173
$ws.Range("C2").Value = 1976

# Row 3 - jose velandia
$ws.Range("A3").Value = "jose velandia"
$ws.Range("B3").Value = 1116458320
$ws.Range("C3").Value = 1996

# Row 4 - oscar velandia
$ws.Range("A4").Value = "oscar velandia"
$ws.Range("B4").Value = 6500471
$ws.Range("C4").Value = 1979

# Column widths (bestFit, matching the authored file). The values below are
# chosen so that, after this engine's internal width rounding, the stored
# <col width="..."> ends up as close as possible to the authored widths
# (17, 11, 10.7109375, 14.140625, 19.42578125, 15.5703125, 20.85546875,
#  13.5703125, 18.85546875, 11.140625).
$ws.Range("A1").EntireColumn.ColumnWidth = 16.166666666666668
$ws.Range("B1").EntireColumn.ColumnWidth = 10.166666666666666
$ws.Range("C1").EntireColumn.ColumnWidth = 9.833333333333332
$ws.Range("D1").EntireColumn.ColumnWidth = 13.333333333333332
$ws.Range("E1").EntireColumn.ColumnWidth = 18.666666666666668
$ws.Range("F1").EntireColumn.ColumnWidth = 14.666666666666666
$ws.Range("G1").EntireColumn.ColumnWidth = 20.0
$ws.Range("H1").EntireColumn.ColumnWidth = 12.666666666666666
$ws.Range("I1").EntireColumn.ColumnWidth = 18.0
$ws.Range("J1").EntireColumn.ColumnWidth = 10.333333333333332

# Selection mirrors the authored file (cell below the data, A5)
$ws.Range("A5").Select()
